$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing (empty) doc_ids column content/formatting for data rows,
# then fill in the new demo doc_ids values for selected rows.
$ws.Range("E2:E43").Clear()

$ws.Range("E6").Value = "statpop-info"
$ws.Range("E7").Value = "pdf-wiki"
$ws.Range("E8").Value = "pdf_online, bevnat-variable"
$ws.Range("E12").Value = "pdf-wiki, tourisme-exemple"
$ws.Range("E14").Value = "bevnat-variable"
$ws.Range("E18").Value = "statpop-info, tourisme-exemple"
$ws.Range("E19").Value = "pop-com-1, pdf_online"
$ws.Range("E22").Value = "pdf-wiki"
$ws.Range("E37").Value = "tourisme-exemple"
$ws.Range("E38").Value = "bevnat-info, tourisme-exemple"
$ws.Range("E40").Value = "statpop-info"

$ws.Range("D4").Select() | Out-Null
